$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Row 9 held "Maquina / carga / recarga la tarjeta ..." and row 10
# was blank. The content now slides down: row 9 becomes the new
# "MaquinaTren / cobroMolineteSubte ..." use case (taking on the
# formatting the blank row 10 had), and row 10 picks up the
# "Maquina / carga ..." entry (taking on row 9's previous
# formatting), with an extra input parameter.
# -----------------------------------------------------------------

# Swap formatting between row 9 and row 10 (using a scratch area so
# row 9's original formatting isn't lost before it is copied down).
$ws.Range("B9:F9").Copy()
$ws.Range("H9:L9").PasteSpecial(-4122)

$ws.Range("B10:F10").Copy()
$ws.Range("B9:F9").PasteSpecial(-4122)

$ws.Range("H9:L9").Copy()
$ws.Range("B10:F10").PasteSpecial(-4122)
$ws.Range("H9:L9").Clear()

# New text for row 9 (MaquinaTren use case).
$ws.Range("B9").Value = "MaquinaTren"
$ws.Range("C9").Value = "cobroMolineteSubte"
$ws.Range("D9").Value = "descuenta el valor del boleto (Subte)"
$ws.Range("E9").Value = "tarjeta: Tarjeta"
$ws.Range("F9").Value = "void"

# New text for row 10 (Maquina use case).
$ws.Range("B10").Value = "Maquina"
$ws.Range("C10").Value = "carga"
$ws.Range("D10").Value = "recarga la tarjeta"
$ws.Range("E10").Value = "tarjeta: Tarjeta, float: valor"
$ws.Range("F10").Value = "void"

# -----------------------------------------------------------------
# Rows 16-18: the standalone "float" return-value entries are
# replaced by the already existing "-" value.
# -----------------------------------------------------------------
$ws.Range("F16").Value = "-"
$ws.Range("F17").Value = "-"
$ws.Range("F18").Value = "-"

# -----------------------------------------------------------------
# Update the active selection to match the edited workbook.
# -----------------------------------------------------------------
$ws.Range("A19").Select()
